$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "vertices"
$ws.Range("B1").Value = "arestas"
$ws.Range("C1").Value = "ponto_x"
$ws.Range("D1").Value = "ponto_y"

# Data rows (vertex id, adjacency list, x, y)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "1,3,4"
$ws.Range("C2").Value = 223
$ws.Range("D2").Value = 157

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "0,2,5"
$ws.Range("C3").Value = 294
$ws.Range("D3").Value = 270

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 356
$ws.Range("D4").Value = 156

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "0,5,6,7"
$ws.Range("C5").Value = 86
$ws.Range("D5").Value = 154

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 297
$ws.Range("D6").Value = 43

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "1,3"
$ws.Range("C7").Value = 146
$ws.Range("D7").Value = 262

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "3,7"
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = 36

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "0,3,6"
$ws.Range("C9").Value = 170
$ws.Range("D9").Value = 41

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = 212
$ws.Range("D10").Value = 351

# Left-align the data block
$ws.Range("A2:D10").HorizontalAlignment = -4131

# Decorative empty column used as a visual gutter (F4:F12), styled with a
# small gray Consolas font, vertically centered
$gutter = $ws.Range("F4:F12")
$gutter.Font.Name = "Consolas"
$gutter.Font.Family = 3
$gutter.Font.Size = 10
$gutter.Font.Color = 13421772
$gutter.VerticalAlignment = -4108

# Column widths (auto best-fit sizing recorded for A and B)
$ws.Columns.Item(1).ColumnWidth = 7.166666666666667
$ws.Columns.Item(2).ColumnWidth = 6.5

# Selection as saved in the workbook
$ws.Range("A1:D10").Select()
